$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated td_sim_1 (C) and record_atd (D) values per corrected relevance markers
# for Appenzeller-Herzog (2019) - van Dis (2020) simulation.

$ws.Range("C2").Value = 12
$ws.Range("D2").Value = 19.5
$ws.Range("C3").Value = 226
$ws.Range("D3").Value = 233
$ws.Range("C4").Value = 109
$ws.Range("D4").Value = 143.5
$ws.Range("C5").Value = 59
$ws.Range("D5").Value = 55
$ws.Range("C6").Value = 47
$ws.Range("D6").Value = 54.5
$ws.Range("C7").Value = 85
$ws.Range("D7").Value = 95
$ws.Range("C8").Value = 106
$ws.Range("D8").Value = 118
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 22
$ws.Range("C10").Value = 187
$ws.Range("D10").Value = 234.5
$ws.Range("C11").Value = 64
$ws.Range("D11").Value = 71
$ws.Range("C12").Value = 68
$ws.Range("D12").Value = 76
$ws.Range("C13").Value = 29
$ws.Range("D13").Value = 36
$ws.Range("C14").Value = 139
$ws.Range("D14").Value = 156
$ws.Range("C15").Value = 95
$ws.Range("D15").Value = 107.5
$ws.Range("C16").Value = 115
$ws.Range("D16").Value = 114.5
$ws.Range("C17").Value = 194
$ws.Range("D17").Value = 197.5
$ws.Range("C18").Value = 15
$ws.Range("D18").Value = 22
$ws.Range("C19").Value = 78
$ws.Range("D19").Value = 86.5
$ws.Range("C20").Value = 34
$ws.Range("D20").Value = 44.5
$ws.Range("C21").Value = 117
$ws.Range("D21").Value = 116
$ws.Range("C22").Value = 852
$ws.Range("D22").Value = 858.5
$ws.Range("C23").Value = 168
$ws.Range("D23").Value = 167.5
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 27
$ws.Range("C25").Value = 178
$ws.Range("D25").Value = 188.5
$ws.Range("C26").Value = 228
$ws.Range("D26").Value = 236
$ws.Range("C27").Value = 88
$ws.Range("D27").Value = 95
$ws.Range("C29").Value = 97
$ws.Range("D29").Value = 110.5
$ws.Range("C30").Value = 24
$ws.Range("D30").Value = 30.5
$ws.Range("C32").Value = 172
$ws.Range("D32").Value = 197
$ws.Range("C33").Value = 62
$ws.Range("D33").Value = 75
$ws.Range("C34").Value = 104
$ws.Range("D34").Value = 116.5
$ws.Range("C35").Value = 60
$ws.Range("D35").Value = 70
$ws.Range("C36").Value = 7
$ws.Range("D36").Value = 30.5
$ws.Range("C37").Value = 6
$ws.Range("D37").Value = 29.5
$ws.Range("C38").Value = 57
$ws.Range("D38").Value = 65.5
$ws.Range("C39").Value = 50
$ws.Range("D39").Value = 51
$ws.Range("C40").Value = 920
$ws.Range("D40").Value = 915
$ws.Range("C41").Value = 71
$ws.Range("D41").Value = 80
$ws.Range("C42").Value = 23
$ws.Range("D42").Value = 30.5
$ws.Range("C43").Value = 53
$ws.Range("D43").Value = 63
$ws.Range("C44").Value = 307
$ws.Range("D44").Value = 331.5
$ws.Range("C45").Value = 73
$ws.Range("D45").Value = 80
$ws.Range("C46").Value = 20
$ws.Range("D46").Value = 26.5
$ws.Range("C47").Value = 265
$ws.Range("D47").Value = 266.5
$ws.Range("C48").Value = 219
$ws.Range("D48").Value = 216
$ws.Range("C49").Value = 113
$ws.Range("D49").Value = 113
$ws.Range("C50").Value = 26
$ws.Range("D50").Value = 33.5
$ws.Range("C51").Value = 320
$ws.Range("D51").Value = 339
$ws.Range("C52").Value = 9
$ws.Range("D52").Value = 16.5
$ws.Range("C53").Value = 21
$ws.Range("D53").Value = 28.5
$ws.Range("C54").Value = 32
$ws.Range("D54").Value = 39
$ws.Range("C55").Value = 196
$ws.Range("D55").Value = 185.5
$ws.Range("C56").Value = 100
$ws.Range("D56").Value = 116
$ws.Range("C57").Value = 195
$ws.Range("D57").Value = 198.5
$ws.Range("C58").Value = 114
$ws.Range("D58").Value = 124.5
$ws.Range("C59").Value = 91
$ws.Range("D59").Value = 107
$ws.Range("C60").Value = 87
$ws.Range("D60").Value = 96.5
$ws.Range("C61").Value = 49
$ws.Range("D61").Value = 59.5
$ws.Range("C62").Value = 279
$ws.Range("D62").Value = 279
$ws.Range("C63").Value = 46
$ws.Range("D63").Value = 59
$ws.Range("C64").Value = 252
$ws.Range("D64").Value = 256.5
$ws.Range("C65").Value = 175
$ws.Range("D65").Value = 182.5
$ws.Range("C66").Value = 83
$ws.Range("D66").Value = 105
$ws.Range("C67").Value = 16
$ws.Range("D67").Value = 23.5
$ws.Range("C68").Value = 45
$ws.Range("D68").Value = 55.5
$ws.Range("C69").Value = 138
$ws.Range("D69").Value = 136
$ws.Range("C70").Value = 52
$ws.Range("D70").Value = 64.5
$ws.Range("C71").Value = 18
$ws.Range("D71").Value = 25
$ws.Range("C72").Value = 99
$ws.Range("D72").Value = 119.5
$ws.Range("C74").Value = 217
$ws.Range("D74").Value = 233
$ws.Range("C75").Value = 65
$ws.Range("D75").Value = 70.5
$ws.Range("C76").Value = 121.7702702702703
